# Add a new customer row (phone 79174401) with 0 points, same shape as the
# existing "no birthday yet" rows (e.g. rows 4-6): phone, blank birthday,
# total_points = 0.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: phone is entered as text (leading apostrophe forces text storage,
# matching the source data which keeps this phone number as a string)
$ws.Range("A8").Value = "'79174401"

# Birthday left blank for this customer
$ws.Range("B8").Value = ""

# Starting points total
$ws.Range("C8").Value = 0
